$wb = $excel.ActiveWorkbook

# --- Sheet "윤다은" (sheet1): fill in row 12 (사용자 개인 페이지구현) and
# row 13 (신청 내역 승인/거절기능 구현) with the new content/result/problem text
# and completion dates, matching the committed worksheet update. ---
$ws1 = $wb.Worksheets.Item("윤다은")

$ws1.Range("B12").Value = "다른 사용자의 개인정보를 열람할 수 있는 페이지를 만듬 -> 예약을 만드는 경우와 수락하는 경우에 참고, 예약된 후 서로 의사소통하기위해 전화번호를 넣었음 -> 채팅기능있으면 필요없음"
$ws1.Range("D12").Value = "2019-06-12"
$ws1.Range("E12").Value = "다른 사용자의 정보 열람 가능"
$ws1.Range("F12").Value = "개인정보 유출의 문제점이 있으므로 채팅을 구현한다면 전화번호 삭제할 예정"
$ws1.Rows.Item(12).RowHeight = 69

$ws1.Range("B13").Value = "공급자에게 온 예약신청을 수락/거절 하는 기능"
$ws1.Range("D13").Value = "2019-06-12"
$ws1.Range("E13").Value = "consumer가 신청한 예약내역을 승인/거절하는 기능 구현"
$ws1.Range("F13").Value = "jsp페이지를 비효율적으로 2개를 만들어서 처리하였음.. 페이지 개수 줄이는 방법이 있을듯"

# --- Sheet "문준범" (sheet2): selection moves from G5 to C4 and loses the
# "tabSelected" flag since the workbook's active tab changes below. ---
$ws2 = $wb.Worksheets.Item("문준범")
$ws2.Activate() | Out-Null
$ws2.Range("C4").Select() | Out-Null

# --- Reactivate "윤다은" as the active sheet/tab, with the new selection
# sitting on H12 (after entering the row-12 values above). ---
$ws1.Activate() | Out-Null
$ws1.Range("H12").Select() | Out-Null
